$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F (dSF) reflecting repulled data
$updates = @{
    3  = 6
    4  = 3
    7  = 4
    10 = -6
    16 = -1
    21 = 1
    22 = 5
    26 = -7
    27 = 10
    28 = -2
    31 = -9
    32 = -2
    38 = -2
    45 = -2
    47 = -5
    54 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
